$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: was scenario17-20230604-003, now scenario17-20230604-001 data
$ws.Range("B7").Value = "CR-PK-CUS-POC-2310041"
$ws.Range("C7").Value = "s1704"
$ws.Range("E7").Value = "scenario1720230604001"
$ws.Range("F7").Value = "PK-CUS-scenario17-20230604-001"
$ws.Range("I7").Value = "PKSUPTOPKCUS17001-s1704-004"
$ws.Range("K7").Value = "PK-SUP-scenario17-20230604-001"

# Row 8: was scenario17-20230604-002, now scenario17-20230604-003 data
$ws.Range("B8").Value = "CR-PK-CUS-POC-2310041"
$ws.Range("C8").Value = "s1704"
$ws.Range("E8").Value = "scenario1720230604003"
$ws.Range("F8").Value = "PK-CUS-scenario17-20230604-003"
$ws.Range("I8").Value = "PKSUPTOPKCUS17001-s1704-004"
$ws.Range("K8").Value = "PK-SUP-scenario17-20230604-003"

# Row 9: was scenario17-20230604-004, now scenario17-20230604-002 data
$ws.Range("B9").Value = "CR-PK-CUS-POC-2310041"
$ws.Range("C9").Value = "s1704"
$ws.Range("E9").Value = "scenario1720230604002"
$ws.Range("F9").Value = "PK-CUS-scenario17-20230604-002"
$ws.Range("I9").Value = "PKSUPTOPKCUS17001-s1704-004"
$ws.Range("K9").Value = "PK-SUP-scenario17-20230604-002"

# Row 10: was scenario17-20230604-001, now scenario17-20230604-004 data
$ws.Range("B10").Value = "CR-PK-CUS-POC-2310041"
$ws.Range("C10").Value = "s1704"
$ws.Range("E10").Value = "scenario1720230604004"
$ws.Range("F10").Value = "PK-CUS-scenario17-20230604-004"
$ws.Range("I10").Value = "PKSUPTOPKCUS17001-s1704-004"
$ws.Range("K10").Value = "PK-SUP-scenario17-20230604-004"
